$d = $word.ActiveDocument

$d.Content.Find.Execute("495÷2=247, 1", $true, $true, $false, $false, $false, $true, 1, $false, "200÷8=25, 0", 2) | Out-Null
$d.Content.Find.Execute("507÷4=126, 3", $true, $true, $false, $false, $false, $true, 1, $false, "173÷8=21, 5", 2) | Out-Null
$d.Content.Find.Execute("496÷2=248, 0", $true, $true, $false, $false, $false, $true, 1, $false, "775÷5=155, 0", 2) | Out-Null
$d.Content.Find.Execute("698÷4=174, 2", $true, $true, $false, $false, $false, $true, 1, $false, "997÷5=199, 2", 2) | Out-Null
$d.Content.Find.Execute("590÷2=295, 0", $true, $true, $false, $false, $false, $true, 1, $false, "509÷4=127, 1", 2) | Out-Null
$d.Content.Find.Execute("707÷3=235, 2", $true, $true, $false, $false, $false, $true, 1, $false, "186÷3=62, 0", 2) | Out-Null
$d.Content.Find.Execute("658÷7=94, 0", $true, $true, $false, $false, $false, $true, 1, $false, "417÷4=104, 1", 2) | Out-Null
$d.Content.Find.Execute("831÷8=103, 7", $true, $true, $false, $false, $false, $true, 1, $false, "175÷7=25, 0", 2) | Out-Null
$d.Content.Find.Execute("111÷6=18, 3", $true, $true, $false, $false, $false, $true, 1, $false, "239÷2=119, 1", 2) | Out-Null
$d.Content.Find.Execute("878÷4=219, 2", $true, $true, $false, $false, $false, $true, 1, $false, "258÷7=36, 6", 2) | Out-Null
$d.Content.Find.Execute("825÷8=103, 1", $true, $true, $false, $false, $false, $true, 1, $false, "728÷5=145, 3", 2) | Out-Null
$d.Content.Find.Execute("152÷9=16, 8", $true, $true, $false, $false, $false, $true, 1, $false, "711÷3=237, 0", 2) | Out-Null
$d.Content.Find.Execute("887÷2=443, 1", $true, $true, $false, $false, $false, $true, 1, $false, "491÷5=98, 1", 2) | Out-Null
$d.Content.Find.Execute("783÷2=391, 1", $true, $true, $false, $false, $false, $true, 1, $false, "732÷3=244, 0", 2) | Out-Null
$d.Content.Find.Execute("590÷7=84, 2", $true, $true, $false, $false, $false, $true, 1, $false, "116÷7=16, 4", 2) | Out-Null
$d.Content.Find.Execute("261÷2=130, 1", $true, $true, $false, $false, $false, $true, 1, $false, "393÷7=56, 1", 2) | Out-Null
$d.Content.Find.Execute("332÷7=47, 3", $true, $true, $false, $false, $false, $true, 1, $false, "261÷5=52, 1", 2) | Out-Null
$d.Content.Find.Execute("219÷2=109, 1", $true, $true, $false, $false, $false, $true, 1, $false, "389÷9=43, 2", 2) | Out-Null
$d.Content.Find.Execute("195÷6=32, 3", $true, $true, $false, $false, $false, $true, 1, $false, "288÷3=96, 0", 2) | Out-Null
$d.Content.Find.Execute("202÷4=50, 2", $true, $true, $false, $false, $false, $true, 1, $false, "905÷6=150, 5", 2) | Out-Null
$d.Content.Find.Execute("464÷2=232, 0", $true, $true, $false, $false, $false, $true, 1, $false, "154÷5=30, 4", 2) | Out-Null
$d.Content.Find.Execute("833÷9=92, 5", $true, $true, $false, $false, $false, $true, 1, $false, "314÷6=52, 2", 2) | Out-Null
$d.Content.Find.Execute("403÷3=134, 1", $true, $true, $false, $false, $false, $true, 1, $false, "845÷4=211, 1", 2) | Out-Null
$d.Content.Find.Execute("728÷2=364, 0", $true, $true, $false, $false, $false, $true, 1, $false, "208÷8=26, 0", 2) | Out-Null
$d.Content.Find.Execute("526÷4=131, 2", $true, $true, $false, $false, $false, $true, 1, $false, "638÷5=127, 3", 2) | Out-Null

Write-Output "Replacements complete"
